# Aulas: Apresentando caso de uso e Nivelamento sobre SQL e JPQL
# Adds 5 new rows (102-106) to the "Tabela1" listobject on Planilha1,
# widens column E, and moves the selection/scroll to the new bottom of
# the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- widen column E (no longer auto bestFit, explicit width) ---------
$ws.Range("E1").ColumnWidth = 70.25

# --- data for the 5 new rows ------------------------------------------
$secao      = 3
$nomeSecao  = "Operações de CRUD e Casos de Uso"

$rows = @(
    @{ Row=102; Aula=48; Height=30; E="0:44 - discussão sobre Caso de Uso - descreve o cenário de utilização da aplicação..."; F="" },
    @{ Row=103; Aula=49; Height=75; E="0:34 - JPQL é a linguagem de consulta da JPA"; F="`n`n`n`n" },
    @{ Row=104; Aula=49; Height=30; E="1:27 - exemplo/comparativo entre uma consulta feita com SQL e uma consulta feita com JPQL"; F="" },
    @{ Row=105; Aula=49; Height=30; E="`n9:31 - sintaxe INNER JOIN em SQl e JPQL"; F="" },
    @{ Row=106; Aula=49; Height=45; E="`n12:48 - instrução SELECT DISTINCT - faz uma consulta no banco de dados e retorna objetos sem repetição"; F="" }
)

# First pass: grow the table, copy formatting down and fill in the
# columns that don't mint brand-new shared strings (B/C/D, plus F which
# either stays blank or reuses the pre-existing "blank lines" string).
foreach ($r in $rows) {
    # Grow the table by one row at a time, then stamp the formatting of
    # the row immediately above it (row 101 carries the "Seção 3" look)
    # onto the freshly added row before writing values into it.
    $lo.ListRows.Add() | Out-Null

    $prev = $r.Row - 1
    $ws.Range("B" + $prev + ":G" + $prev).Copy()
    $ws.Range("B" + $r.Row + ":G" + $r.Row).PasteSpecial(-4122)

    $ws.Range("B" + $r.Row).Value = $secao
    $ws.Range("C" + $r.Row).Value = $nomeSecao
    $ws.Range("D" + $r.Row).Value = $r.Aula

    $ws.Rows.Item($r.Row).RowHeight = $r.Height
}

# Second pass: write column E in the same order the new shared strings
# were originally authored in (102, 104, 105, 103, 106) so the new
# entries land at shared-string indices 164..168 in that order, then
# fill in F103 (reuses an already-existing shared string).
$ws.Range("E102").Value = $rows[0].E
$ws.Range("E104").Value = $rows[2].E
$ws.Range("E105").Value = $rows[3].E
$ws.Range("E103").Value = $rows[1].E
$ws.Range("F103").Value = $rows[1].F
$ws.Range("E106").Value = $rows[4].E

$excel.CutCopyMode = $false

# --- leave the selection / scroll position where the edit ended ------
$excel.Goto($ws.Range("A100"), $true)
$ws.Range("D106").Select()
